# Apply the two changes described by the diff:
#  1. Update the cached "datetimeFigureOut" auto-date field text from
#     12/16/2022 -> 12/17/2022 everywhere it is cached (slide master +
#     all 5 slide layouts).
#  2. Delete the "object 10" text box (shape id 10) from slide 4 - the
#     paragraph describing "Operacao: Termino de Obra ... LTV de 41%."

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "12/16/2022") {
                $tr.Text = "12/17/2022"
            }
        }
    }
}

# 1a. Slide master holder that carries the cached date field.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# 1b. Every slide layout under the master also caches the same field.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# 2. Remove the "object 10" textbox from slide 4.
$slide4 = $p.Slides.Item(4)
for ($i = $slide4.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide4.Shapes.Item($i)
    if ($sh.Name -eq "object 10") {
        $sh.Delete()
    }
}
